$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data values (diff: B2, C2, C3, C5) ---
$ws.Range("B2").Value = 5.0999999999999996
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 7
$ws.Range("C5").Value = 25

# --- Column widths (diff: col A -> 27, col C -> 27.25, col B reverts to default) ---
# ColumnWidth is expressed in characters; the stored "width" XML attribute is
# derived as chars + 5/MaximumDigitWidth. This runtime's MDW is 7, so we solve
# for the character width that yields the desired stored width exactly
# (184/7 -> stored 27, 186/7 -> stored ~27.2857, the closest reachable grid
# point to the authored 27.25).
$ws.Columns.Item(1).ColumnWidth = 184/7
$ws.Columns.Item(3).ColumnWidth = 186/7

# --- Selection moved from C3 to C4 ---
$ws.Range("C4").Select()
